# Generate Report for Handoff
# Adds two new handed-off files (da4cdb8d... and e863d0a5...) as new rows
# to the Overview, zh-cn and de-de tables.

$wb = $excel.ActiveWorkbook

$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

$loOv = $ov.ListObjects.Item(1)
$loZh = $zh.ListObjects.Item(1)
$loDe = $de.ListObjects.Item(1)

# Hyperlink colour used throughout this workbook (Cornflower Blue, FF6495ED)
$linkColor = 15570276

function Style-AsLink($rng) {
    $rng.Font.Underline = 2
    $rng.Font.Color = $linkColor
}

# ---------------------------------------------------------------------------
# Overview sheet — 2 new rows (row 4 + row 5)
# Columns: A File Name | B Path And Name | C Extension | D Publish URL |
#          E zh-cn | F de-de | G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$loOv.ListRows.Add() | Out-Null
$ov.Range("A4").Value = "da4cdb8d-2d43-4fba-9e04-d896cb4a9270.md"
$ov.Range("B4").Value = "e2e\da4cdb8d-2d43-4fba-9e04-d896cb4a9270.md"
$ov.Range("C4").Value = ".md"
$ov.Range("E4").Value = "Ready for handoff"
$ov.Range("F4").Value = "Ready for handoff"
$ov.Range("G4").Value = "2016-09-01 02:49:25"
$ov.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ov.Hyperlinks.Add($ov.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec55f1e4c7c1cbd4e938b3fc0aeeeb5fb9af874d/e2e/da4cdb8d-2d43-4fba-9e04-d896cb4a9270.md", "", "", "e2e\da4cdb8d-2d43-4fba-9e04-d896cb4a9270.md") | Out-Null
Style-AsLink $ov.Range("B4")

$loOv.ListRows.Add() | Out-Null
$ov.Range("A5").Value = "e863d0a5-2dc3-4fa6-a0fa-4b11f1b92619.md"
$ov.Range("B5").Value = "e2e\e863d0a5-2dc3-4fa6-a0fa-4b11f1b92619.md"
$ov.Range("C5").Value = ".md"
$ov.Range("E5").Value = "Ready for handoff"
$ov.Range("F5").Value = "Ready for handoff"
$ov.Range("G5").Value = "2016-09-01 02:49:25"
$ov.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ov.Hyperlinks.Add($ov.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec55f1e4c7c1cbd4e938b3fc0aeeeb5fb9af874d/e2e/e863d0a5-2dc3-4fa6-a0fa-4b11f1b92619.md", "", "", "e2e\e863d0a5-2dc3-4fa6-a0fa-4b11f1b92619.md") | Out-Null
Style-AsLink $ov.Range("B5")

# ---------------------------------------------------------------------------
# zh-cn sheet — 2 new rows (row 4 + row 5)
# Columns: A Source File Name | B File Extension | C Status | D Source Path |
#          E Priority | F Content Duplicate | G Latest Handoff File |
#          H Latest Handoff Datetime | I Latest Target File |
#          J Latest Handback File | K Latest Handback DateTime |
#          L Reference Tokens | M To be localized | N Dependency From |
#          O Has metadata | P Error Detail
# ---------------------------------------------------------------------------
$loZh.ListRows.Add() | Out-Null
$zh.Range("A4").Value = "da4cdb8d-2d43-4fba-9e04-d896cb4a9270.md"
$zh.Range("B4").Value = ".md"
$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("D4").Value = "e2e"
$zh.Range("E4").Value = "ht"
$zh.Range("F4").Value = "False"
$zh.Range("G4").Value = "da4cdb8d-2d43-4fba-9e04-d896cb4a9270.6b2b737368d5d4c5755d5d3d8e0f48b4de47b0af.zh-cn.xlf"
$zh.Range("H4").Value = "2016-09-01 02:49:20"
$zh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("I4").Value = ""
$zh.Range("J4").Value = ""
$zh.Range("K4").Value = "0001-01-01 00:00:00"
$zh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("L4").Value = ""
$zh.Range("M4").Value = "True"
$zh.Range("N4").Value = ""
$zh.Range("O4").Value = "False"
$zh.Range("P4").Value = ""
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec55f1e4c7c1cbd4e938b3fc0aeeeb5fb9af874d/e2e/da4cdb8d-2d43-4fba-9e04-d896cb4a9270.md", "", "", "da4cdb8d-2d43-4fba-9e04-d896cb4a9270.md") | Out-Null
Style-AsLink $zh.Range("A4")

$loZh.ListRows.Add() | Out-Null
$zh.Range("A5").Value = "e863d0a5-2dc3-4fa6-a0fa-4b11f1b92619.md"
$zh.Range("B5").Value = ".md"
$zh.Range("C5").Value = "Ready for handoff"
$zh.Range("D5").Value = "e2e"
$zh.Range("E5").Value = "ht"
$zh.Range("F5").Value = "False"
$zh.Range("G5").Value = "e863d0a5-2dc3-4fa6-a0fa-4b11f1b92619.8ecfd519672ae8fcdd0d581be0ad55d76debe0a2.zh-cn.xlf"
$zh.Range("H5").Value = "2016-09-01 02:49:20"
$zh.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("I5").Value = ""
$zh.Range("J5").Value = ""
$zh.Range("K5").Value = "0001-01-01 00:00:00"
$zh.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("L5").Value = ""
$zh.Range("M5").Value = "True"
$zh.Range("N5").Value = ""
$zh.Range("O5").Value = "False"
$zh.Range("P5").Value = ""
$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec55f1e4c7c1cbd4e938b3fc0aeeeb5fb9af874d/e2e/e863d0a5-2dc3-4fa6-a0fa-4b11f1b92619.md", "", "", "e863d0a5-2dc3-4fa6-a0fa-4b11f1b92619.md") | Out-Null
Style-AsLink $zh.Range("A5")

# ---------------------------------------------------------------------------
# de-de sheet — 2 new rows (row 4 + row 5), same column layout as zh-cn
# ---------------------------------------------------------------------------
$loDe.ListRows.Add() | Out-Null
$de.Range("A4").Value = "da4cdb8d-2d43-4fba-9e04-d896cb4a9270.md"
$de.Range("B4").Value = ".md"
$de.Range("C4").Value = "Ready for handoff"
$de.Range("D4").Value = "e2e"
$de.Range("E4").Value = "ht"
$de.Range("F4").Value = "False"
$de.Range("G4").Value = "da4cdb8d-2d43-4fba-9e04-d896cb4a9270.6b2b737368d5d4c5755d5d3d8e0f48b4de47b0af.de-de.xlf"
$de.Range("H4").Value = "2016-09-01 02:49:25"
$de.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("I4").Value = ""
$de.Range("J4").Value = ""
$de.Range("K4").Value = "0001-01-01 00:00:00"
$de.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("L4").Value = ""
$de.Range("M4").Value = "True"
$de.Range("N4").Value = ""
$de.Range("O4").Value = "False"
$de.Range("P4").Value = ""
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec55f1e4c7c1cbd4e938b3fc0aeeeb5fb9af874d/e2e/da4cdb8d-2d43-4fba-9e04-d896cb4a9270.md", "", "", "da4cdb8d-2d43-4fba-9e04-d896cb4a9270.md") | Out-Null
Style-AsLink $de.Range("A4")

$loDe.ListRows.Add() | Out-Null
$de.Range("A5").Value = "e863d0a5-2dc3-4fa6-a0fa-4b11f1b92619.md"
$de.Range("B5").Value = ".md"
$de.Range("C5").Value = "Ready for handoff"
$de.Range("D5").Value = "e2e"
$de.Range("E5").Value = "ht"
$de.Range("F5").Value = "False"
$de.Range("G5").Value = "e863d0a5-2dc3-4fa6-a0fa-4b11f1b92619.8ecfd519672ae8fcdd0d581be0ad55d76debe0a2.de-de.xlf"
$de.Range("H5").Value = "2016-09-01 02:49:25"
$de.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("I5").Value = ""
$de.Range("J5").Value = ""
$de.Range("K5").Value = "0001-01-01 00:00:00"
$de.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("L5").Value = ""
$de.Range("M5").Value = "True"
$de.Range("N5").Value = ""
$de.Range("O5").Value = "False"
$de.Range("P5").Value = ""
$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec55f1e4c7c1cbd4e938b3fc0aeeeb5fb9af874d/e2e/e863d0a5-2dc3-4fa6-a0fa-4b11f1b92619.md", "", "", "e863d0a5-2dc3-4fa6-a0fa-4b11f1b92619.md") | Out-Null
Style-AsLink $de.Range("A5")

Write-Host "Done adding handback rows"
